# Auto-generated edit script: updates crypto price/volume table cells
# to match the new scraped values, keeping every cell stored as text
# (matches the original inlineStr/shared-string text cells; avoids Excel
# auto-converting numeric-looking strings like "211.39" or "8.037" into
# real numbers, which would lose formatting such as trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextCell "D2" "26.209.42"
Set-TextCell "E2" "  -0.84%  "
Set-TextCell "D3" "1.680.93"
Set-TextCell "E3" "  -0.68%  "
Set-TextCell "E4" "  -0.60%  "
Set-TextCell "D5" "211.39"
Set-TextCell "E5" "  -3.44%  "
Set-TextCell "D6" "0.5307"
Set-TextCell "E6" "  -3.93%  "
Set-TextCell "E7" "  -0.61%  "
Set-TextCell "D8" "0.2688"
Set-TextCell "E8" "  -1.41%  "
Set-TextCell "D9" "0.06311"
Set-TextCell "E9" "  -2.66%  "
Set-TextCell "D10" "21.36"
Set-TextCell "E10" "  -3.50%  "
Set-TextCell "D11" "0.07535"
Set-TextCell "E11" "  -2.19%  "
Set-TextCell "D12" "1.680.05"
Set-TextCell "E12" "  -0.66%  "
Set-TextCell "D13" "4.486"
Set-TextCell "E13" "  -1.46%  "
Set-TextCell "D14" "0.5675"
Set-TextCell "E14" "  -2.65%  "
Set-TextCell "D15" "0.000008145"
Set-TextCell "E15" "  -3.49%  "
Set-TextCell "D16" "66.46"
Set-TextCell "E16" "  +1.79%  "
Set-TextCell "D17" "26.260.85"
Set-TextCell "E17" "  -0.86%  "
Set-TextCell "E18" "  -0.62%  "
Set-TextCell "D19" "4.860"
Set-TextCell "E19" "  -1.91%  "
Set-TextCell "D20" "10.58"
Set-TextCell "E20" "  -3.86%  "
Set-TextCell "D21" "189.07"
Set-TextCell "E21" "  -0.64%  "
Set-TextCell "D22" "6.227"
Set-TextCell "E22" "  -0.25%  "
Set-TextCell "D23" "1.006"
Set-TextCell "E23" "  -0.60%  "
Set-TextCell "D24" "147.61"
Set-TextCell "E24" "  -1.66%  "
Set-TextCell "D25" "0.1265"
Set-TextCell "E25" "  -3.34%  "
Set-TextCell "D26" "7.625"
Set-TextCell "E26" "  -3.54%  "
Set-TextCell "D27" "15.90"
Set-TextCell "E27" "  +0.96%  "
Set-TextCell "D28" "0.06458"
Set-TextCell "E28" "  +1.92%  "
Set-TextCell "D29" "1.343"
Set-TextCell "E29" "  -6.00%  "
Set-TextCell "D30" "1.287"
Set-TextCell "E30" "  -3.35%  "
Set-TextCell "D31" "3.540"
Set-TextCell "E31" "  -1.49%  "
Set-TextCell "D32" "3.488"
Set-TextCell "E32" "  -2.88%  "
Set-TextCell "D33" "1.657"
Set-TextCell "E33" "  -1.26%  "
Set-TextCell "D34" "1.012"
Set-TextCell "E34" "  -3.04%  "
Set-TextCell "D35" "0.6130"
Set-TextCell "E35" "  -1.35%  "
Set-TextCell "E36" "  +0.44%  "
Set-TextCell "D37" "2.718"
Set-TextCell "E37" "  -0.27%  "
Set-TextCell "D38" "6.188"
Set-TextCell "E38" "  -0.74%  "
Set-TextCell "D39" "0.01622"
Set-TextCell "E39" "  -1.19%  "
Set-TextCell "D40" "1.103.66"
Set-TextCell "E40" "  -1.87%  "
Set-TextCell "D41" "0.8678"
Set-TextCell "E41" "  -1.59%  "
Set-TextCell "E42" "  -0.98%  "
Set-TextCell "E43" "  -0.63%  "
Set-TextCell "D44" "1.834.00"
Set-TextCell "E44" "  -0.51%  "
Set-TextCell "D45" "0.00000000108"
Set-TextCell "E45" "  +0.08%  "
Set-TextCell "D46" "57.22"
Set-TextCell "E46" "  -0.43%  "
Set-TextCell "D47" "1.004"
Set-TextCell "E47" "  -0.46%  "
Set-TextCell "B48" "EnergySwap"
Set-TextCell "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D48" "8.037"
Set-TextCell "E48" "  -2.36%  "
Set-TextCell "B49" "Cronos"
Set-TextCell "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D49" "0.05264"
Set-TextCell "E49" "  -0.38%  "
Set-TextCell "D50" "0.4273"
Set-TextCell "E50" "  -1.82%  "
